# Fixed a bug in respin: correct the reel-strip frequency table so that
# each row's data (symbol id + per-reel counts) is restored to the
# correct symbol ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(501, 9, 52, 30, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(1203, 3, 15, 15, 15, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(401, 9, 48, 67, 75, 45)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}

$data2 = @(
    @(3, 0, 3, 3, 3, 3),
    @(1101, 0, 15, 30, 30, 0),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(1, 0, 2, 2, 2, 2)
)

$row = 16
foreach ($r in $data2) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row = $row + 1
}
